$wb = $excel.ActiveWorkbook

# Sheets "展览" (exhibition) and "全部类型" (all types) both list the same
# two events in rows 2-3; column F ("想去人数" / "want to go" count) went up
# for each of them.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 154
    $ws.Range("F3").Value = 106
}
